## 29-ABR - 17:43 - JCSantos
## Updates Settings values (new Gmail/ChatGPT credentials), inserts a new
## "ChatGPT"/"CredentialCHATGPT" settings row, refreshes a couple of row
## heights (wrapped-text rows) and leaves the workbook with the "Assets"
## sheet active/selected, matching the author's last on-screen state.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Settings")
$ws2 = $wb.Worksheets.Item("Constants")
$ws3 = $wb.Worksheets.Item("Assets")

# --- Settings sheet: update existing values -------------------------------
$ws1.Range("B2").Value = "Recomendaciones"
$ws1.Range("B2").Style = "Normal"
$ws1.Range("B3").Value = "nexusdatacorporation@gmail.com's workspace"
$ws1.Range("B9").Value = "CredentialGmail"

# --- Settings sheet: insert a new row for the ChatGPT credential ----------
# (pushes the old "ListKill" / "EXCEL" row from 10 down to 11)
$ws1.Rows.Item(10).Insert()
$ws1.Rows.Item(10).RowHeight = 14.25
$ws1.Range("A10").Value = "ChatGPT"
$ws1.Range("B10").Value = "CredentialCHATGPT"

# --- Row heights for wrapped-text description rows ------------------------
$ws1.Rows.Item(3).RowHeight = 45
$ws1.Rows.Item(5).RowHeight = 30

$ws2.Rows.Item(2).RowHeight = 30
$ws2.Rows.Item(3).RowHeight = 45
$ws2.Rows.Item(17).RowHeight = 45

# --- Selections / active sheet --------------------------------------------
$ws1.Range("B10").Select()
$ws2.Range("C7").Select()
$ws3.Range("A2:D3").Select()
